$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..42 down to 4..43
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the weekly price record (same template as the
# surrounding rows, new date/quality/volume/price values)
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "Macroferia Regional de Talca"
$ws.Range("C3").Value = "Maule"
$ws.Range("D3").Value = 44490
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 25000
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 2500
$ws.Range("T3").Value = 10
